$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item("TextBox 67")
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Locate the unique sentence "GPIO adaptor extract info and output to LED, Servo or "
# so the offsets below are resolved robustly instead of relying on a hard-coded
# paragraph index.
$anchor = $tr.Find("GPIO adaptor extract")
$adaptorStart = $anchor.Start + 5

# Original run text: "GPIO " + "adaptor" + " extract info and output to LED, Servo or "
# Target text:       "GPIO " + "plugin module" + " " + "extract info and output to LED, Servo or "

# Replace "adaptor" with "plugin module"
$word = $tr.Characters($adaptorStart, 7)
$word.Text = "plugin module"

# Split the following space into its own run, so the final run layout is:
#   "GPIO " / "plugin module" / " " / "extract info and output to LED, Servo or "
$space = $tr.Characters($adaptorStart + 13, 1)
$space.Text = " "
